# Generate Report for Handoff
# - Flip the "In Translation" status to "Ready for handoff" (Overview + per-locale sheets)
# - Refresh the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
# - Widen the now-longer "Status" columns so the new text fits

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-04 05:05:29"
# Raw stored width target is 17.2159881591797 characters; the COM ColumnWidth
# setter quantizes to the nearest 1/6-character increment, so feed it the
# input that lands closest to that stored value.
$overview.Range("E:E").ColumnWidth = 16.333333333333336
$overview.Range("F:F").ColumnWidth = 16.333333333333336

# --- zh-cn sheet ------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-04 05:05:25"
$zhcn.Range("C:C").ColumnWidth = 16.333333333333336

# --- de-de sheet ------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-04 05:05:29"
$dede.Range("C:C").ColumnWidth = 16.333333333333336
